# Insert two new header columns ("xm" and "abo_i") before the existing
# "adequacy" column (currently column G), shifting all subsequent columns
# to the right by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at G:H, pushing existing G:H (and beyond) to the right.
# Using EntireColumn.Insert twice so the new cells inherit formatting from
# the columns they are inserted in front of (standard Excel behaviour).
$ws.Range("G:H").EntireColumn.Insert()

# Set the new header values in row 1.
$ws.Range("G1").Value = "xm"
$ws.Range("H1").Value = "abo_i"
